# Commit: "change template in french"
# Update the English-only column headers in row 1 of Sheet1 to bilingual
# English/French headers. Cell positions and styles are unchanged; only
# the header text values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Children/Enfants (5-17)"
$ws.Range("D1").Value = "Girls/Filles  (5-17)"
$ws.Range("E1").Value = "Boys/Garcons  (5-17)"
$ws.Range("F1").Value = "Host/Hôte-- Children (5-17)"
$ws.Range("G1").Value = "Host/Hôte -- Girls  (5-17)"
$ws.Range("H1").Value = "Host/Hôte -- Boys  (5-17)"
$ws.Range("I1").Value = "IDP/PDI -- Children (5-17)"
$ws.Range("J1").Value = "IDP/PDI -- Girls  (5-17)"
$ws.Range("K1").Value = "IDP/PDI -- Boys  (5-17)"
$ws.Range("L1").Value = "Returnees/Retournés -- Children (5-17)"
$ws.Range("M1").Value = "Returnees/Retournés -- Girls  (5-17)"
$ws.Range("N1").Value = "Returnees/Retournés -- Boys  (5-17)"
$ws.Range("O1").Value = "Refugees/Refugiee -- Children (5-17)"
$ws.Range("P1").Value = "Refugees/Refugiee -- Girls  (5-17)"
$ws.Range("Q1").Value = "Refugees/Refugiee -- Boys  (5-17)"

# Match the author's resulting cell selection/active-cell state.
$ws.Range("D16").Select()
